$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 104
$ws.Cells.Item(3, 1).Value = 106
$ws.Cells.Item(4, 1).Value = 108
$ws.Cells.Item(5, 1).Value = 110
$ws.Cells.Item(6, 1).Value = 111
$ws.Cells.Item(7, 1).Value = 113
$ws.Cells.Item(8, 1).Value = 116
$ws.Cells.Item(9, 1).Value = 119
$ws.Cells.Item(10, 1).Value = 122
$ws.Cells.Item(11, 1).Value = 124
$ws.Cells.Item(12, 1).Value = 126
$ws.Cells.Item(13, 1).Value = 128
$ws.Cells.Item(14, 1).Value = 21
$ws.Cells.Item(15, 1).Value = 36
$ws.Cells.Item(16, 1).Value = 91
$ws.Cells.Item(17, 1).Value = 170
$ws.Cells.Item(18, 1).Value = 185
$ws.Cells.Item(19, 1).Value = 224
$ws.Cells.Item(20, 1).Value = 277
$ws.Cells.Item(21, 1).Value = 287
$ws.Cells.Item(22, 1).Value = 363
$ws.Cells.Item(23, 1).Value = 375
$ws.Cells.Item(24, 1).Value = 395
$ws.Cells.Item(25, 1).Value = 424
$ws.Cells.Item(26, 1).Value = 465
$ws.Cells.Item(27, 1).Value = 487
